$d = $word.ActiveDocument

# --- 1. Text corrections: drop the stray double-comma in these row labels ---
$d.Content.Find.Execute("Ramets,, before flowering", $false, $false, $false, $false, $false, $true, 1, $false, "Ramets before flowering", 2) | Out-Null
$d.Content.Find.Execute("Ramets,, after flowering", $false, $false, $false, $false, $false, $true, 1, $false, "Ramets after flowering", 2) | Out-Null
$d.Content.Find.Execute("Herbivory, before flowering (binary)", $false, $false, $false, $false, $false, $true, 1, $false, "Herbivory before flowering (binary)", 2) | Out-Null
$d.Content.Find.Execute("Herbivory, before flowering (quantitative)", $false, $false, $false, $false, $false, $true, 1, $false, "Herbivory before flowering (quantitative)", 2) | Out-Null
$d.Content.Find.Execute("Herbivory, after flowering (binary)", $false, $false, $false, $false, $false, $true, 1, $false, "Herbivory after flowering (binary)", 2) | Out-Null
$d.Content.Find.Execute("Herbivory, after flowering (quantitative)", $false, $false, $false, $false, $false, $true, 1, $false, "Herbivory after flowering (quantitative)", 2) | Out-Null

# --- 2. Table column widths: first column narrows, last column widens ---
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 4205 / 20.0
$t.Columns.Item(7).Width = 977 / 20.0

# --- 3. New results row: "Pollinaria removed" ---
$t.Rows.Add() | Out-Null
$i = $t.Rows.Count

$t.Cell($i, 1).Range.Text = "Pollinaria removed"
$t.Cell($i, 2).Range.Text = "0.017"
$t.Cell($i, 3).Range.Text = "0.449"
$t.Cell($i, 4).Range.Text = "1.003"
$t.Cell($i, 5).Range.Text = "3.304"

$c6 = $t.Cell($i, 6)
$c6.Range.Text = "0.0345"
$boldRng = $d.Range($c6.Range.Start, $c6.Range.End - 1)
$boldRng.Font.Bold = 1

$t.Cell($i, 7).Range.Text = "27.192"

Write-Host "Edit complete"
